$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete row 290 (MFE / McAllen, TX entry) entirely; all rows below shift up by one.
$ws.Rows.Item(290).Delete()
